# semana 40 de 2025
# Adds week 39 (column AP) and week 40 (column AQ) to the weekly
# surveillance sheet, mirroring the existing week columns (1..38 in
# D1:AO1) and extending the per-row counts that already reach column AO.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new week numbers, stored as text (matching the
# existing week-number headers in D1:AO1) with the same bold/centered
# style (s="1") as the rest of the header row.
# The leading apostrophe forces text storage instead of Excel's automatic
# numeric coercion; pasting AO1's format on top afterwards normalizes the
# style back onto the shared header style instead of a one-off variant.
$ws.Cells.Item(1, 42).Value = "'39"   # AP1 -> text "39"
$ws.Cells.Item(1, 43).Value = "'40"   # AQ1 -> text "40"
$ws.Range("AO1").Copy() | Out-Null
$ws.Range("AP1:AQ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows: AP (week 39) and AQ (week 40) counts.
# Only rows whose data previously extended through column AO (the last
# existing week column) gained new values; a handful only received one
# of the two new columns, matching the source data exactly.
$weekData = @{
    2  = @{ AP = 0;  AQ = 0 }
    3  = @{ AP = 0;  AQ = 0 }
    5  = @{ AP = 0;  AQ = 0 }
    6  = @{ AP = 2;  AQ = 2 }
    7  = @{ AP = 0;  AQ = 0 }
    8  = @{ AP = 0;  AQ = 0 }
    9  = @{ AP = 0;  AQ = 0 }
    10 = @{ AP = 0;  AQ = 0 }
    11 = @{ AP = 0 }
    12 = @{ AP = 0;  AQ = 0 }
    13 = @{ AP = 0 }
    14 = @{ AP = 0;  AQ = 0 }
    15 = @{ AQ = 0 }
    16 = @{ AP = 0;  AQ = 0 }
    17 = @{ AP = 0;  AQ = 0 }
    22 = @{ AP = 0;  AQ = 0 }
    23 = @{ AP = 0;  AQ = 0 }
    24 = @{ AQ = 0 }
    25 = @{ AP = 0;  AQ = 0 }
    26 = @{ AP = 0 }
    28 = @{ AP = 52; AQ = 0 }
    29 = @{ AP = 1;  AQ = 0 }
    30 = @{ AP = 0;  AQ = 1 }
    31 = @{ AP = 0;  AQ = 0 }
    35 = @{ AP = 1;  AQ = 4 }
    36 = @{ AP = 0;  AQ = 0 }
    37 = @{ AP = 0;  AQ = 0 }
    38 = @{ AP = 0;  AQ = 0 }
    41 = @{ AP = 0;  AQ = 0 }
    42 = @{ AP = 0;  AQ = 0 }
    43 = @{ AP = 0;  AQ = 0 }
    44 = @{ AP = 0 }
    45 = @{ AP = 0;  AQ = 0 }
    46 = @{ AP = 0;  AQ = 0 }
    47 = @{ AP = 0;  AQ = 0 }
    48 = @{ AP = 0;  AQ = 0 }
    49 = @{ AP = 0;  AQ = 0 }
    50 = @{ AP = 0;  AQ = 0 }
    51 = @{ AP = 0;  AQ = 0 }
    52 = @{ AP = 0;  AQ = 0 }
    53 = @{ AP = 0;  AQ = 0 }
    54 = @{ AP = 0;  AQ = 0 }
    55 = @{ AP = 0;  AQ = 0 }
    56 = @{ AP = 0;  AQ = 0 }
    57 = @{ AP = 0;  AQ = 0 }
    58 = @{ AP = 0;  AQ = 0 }
}

foreach ($row in $weekData.Keys) {
    $cols = $weekData[$row]
    if ($cols.ContainsKey("AP")) {
        $ws.Cells.Item($row, 42).Value = $cols["AP"]
    }
    if ($cols.ContainsKey("AQ")) {
        $ws.Cells.Item($row, 43).Value = $cols["AQ"]
    }
}

# --- Extend the sheet's used dimension to include the new columns.
$ws.Range("A1:AQ58").Select() | Out-Null
